$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.569.70"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.31%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.017.29"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.58%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "585.13"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.54%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.18"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.89%  "

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.522"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.95%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.016.08"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.57%  "

$ws.Range("E10").Value = "  -3.72%  "

$ws.Range("E11").Value = "  -1.19%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.442"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.68%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000230"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.89%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.86"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -5.54%  "

$ws.Range("E15").Value = "  +2.51%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.505.31"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.90%  "

$ws.Range("B17").Value = "Polkadot"
$ws.Range("C17").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.03"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.17%  "

$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "62.476.82"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.54%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.010.87"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.79%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "465.42"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.49%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.95"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.47%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.686"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.54%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.36"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.92%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.34"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.13%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "79.93"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.87%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.40"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.17%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.21"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.78%  "

$ws.Range("E28").Value = "  +0.07%  "

$ws.Range("E29").Value = "  +0.02%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.63"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.96%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.16"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.42%  "

$ws.Range("E32").Value = "  -0.97%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.49"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.46%  "

$ws.Range("E34").Value = "  -4.59%  "

$ws.Range("E35").Value = "  -1.02%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0₃0797"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.76%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.76"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.74%  "

$ws.Range("E38").Value = "  -3.13%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "50.40"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.24%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.96"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.63%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.94"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -10.71%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "421.93"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.25%  "

$ws.Range("E43").Value = "  +0.08%  "

$ws.Range("E44").Value = "  -2.81%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.780.83"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.43%  "

$ws.Range("E46").Value = "  -1.73%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "38.31"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.09%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "129.27"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.46%  "

$ws.Range("E49").Value = "  +0.03%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "24.16"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.83%  "

$ws.Range("E51").Value = "  -1.02%  "
